$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 104
$ws.Range("F6").Value = 1869
$ws.Range("F7").Value = 865
$ws.Range("F8").Value = 1356
$ws.Range("F10").Value = 771
$ws.Range("F12").Value = 2915
$ws.Range("F13").Value = 387
$ws.Range("F14").Value = 880
$ws.Range("F15").Value = 1132
$ws.Range("F17").Value = 335
$ws.Range("F19").Value = 1664
$ws.Range("F20").Value = 344
$ws.Range("F21").Value = 1266
$ws.Range("F22").Value = 213
$ws.Range("F25").Value = 1073
$ws.Range("F26").Value = 1527
$ws.Range("F27").Value = 1474
$ws.Range("F29").Value = 346
$ws.Range("F30").Value = 1299
$ws.Range("F31").Value = 448
$ws.Range("F35").Value = 1856
$ws.Range("F36").Value = 485
$ws.Range("F40").Value = 2300
$ws.Range("F41").Value = 152
$ws.Range("F43").Value = 2804
$ws.Range("F46").Value = 646
$ws = $wb.Worksheets.Item(2)
$ws.Range("F5").Value = 63
$ws.Range("F7").Value = 29
$ws.Range("F13").Value = 113808
$ws.Range("F17").Value = 73
$ws.Range("F18").Value = 73
$ws.Range("F20").Value = 290
$ws.Range("F22").Value = 285
$ws.Range("F25").Value = 71
$ws.Range("F26").Value = 66
$ws.Range("F27").Value = 66
$ws.Range("F30").Value = 47
$ws.Range("F31").Value = 133
$ws.Range("F36").Value = 89
$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 297
$ws.Range("F6").Value = 4859
$ws.Range("F9").Value = 681
$ws.Range("F10").Value = 951
$ws.Range("F11").Value = 551
$ws.Range("F12").Value = 642
$ws.Range("F13").Value = 1370
$ws.Range("F14").Value = 389
$ws.Range("F15").Value = 1280
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 297
$ws.Range("F4").Value = 4859
$ws.Range("F5").Value = 681
$ws.Range("F6").Value = 951
$ws.Range("F7").Value = 551
$ws.Range("F8").Value = 104
$ws.Range("F9").Value = 642
$ws.Range("F10").Value = 1370
$ws.Range("F11").Value = 1869
$ws.Range("F12").Value = 865
$ws.Range("F13").Value = 1356
$ws.Range("F14").Value = 771
$ws.Range("F15").Value = 771
$ws.Range("F16").Value = 1280
$ws.Range("F17").Value = 2915
$ws.Range("F19").Value = 387
$ws.Range("F20").Value = 880
$ws.Range("F21").Value = 1132
$ws.Range("F23").Value = 335
$ws.Range("F24").Value = 1664
$ws.Range("F26").Value = 344
$ws.Range("F28").Value = 1266
$ws.Range("F29").Value = 213
$ws.Range("F32").Value = 1527
$ws.Range("F33").Value = 1474
$ws.Range("F35").Value = 346
$ws.Range("F36").Value = 73
$ws.Range("F37").Value = 1299
$ws.Range("F38").Value = 448
$ws.Range("F41").Value = 1856
$ws.Range("F42").Value = 66
$ws.Range("F43").Value = 133
$ws.Range("F45").Value = 2300
$ws.Range("F47").Value = 2804
$ws.Range("F49").Value = 646
